$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.281.41"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "1.864.53"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.84"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4720"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2910"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06556"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.91"
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07931"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.05"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "1.869.05"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.160"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6828"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.05"
$ws.Range("E16").Value = "  -5.12%  "
$ws.Range("D17").Value = "30.257.43"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.75"
$ws.Range("E18").Value = "  +8.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007422"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "2.112.30"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.310"
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.191"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.55"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.239"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.95"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.956"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09861"
$ws.Range("E30").Value = "  +1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.381"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.473"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04718"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7040"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01882"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.612"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.46"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.950"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8449"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4169"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.45"
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.186"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "950.98"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.216"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.18"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05659"
$ws.Range("E51").Value = "  +0.39%  "
